$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B9").Value = "If he had pursued his passion, he might be happier now."
$ws.Range("C9").Value = "Se tivesse seguido sua paixão, poderia estar mais feliz agora."

$ws.Range("B10").Value = "She wishes she had traveled more before starting a family."
$ws.Range("C10").Value = "Ela gostaria de ter viajado mais antes de começar uma família."
